# Increase number of template sheets by 1
#
# Duplicate the "WFA (6)" sheet (the last WFA-style per-player worksheet,
# header rows + "Table145678" list object + conditional formatting on I4)
# to create a new "WFA (7)" sheet, inserted right before "Player Summary" -
# the same thing Excel's "Move or Copy... > Create a copy" does when you
# drag a new copy of the last WFA tab into place.

$wb = $excel.ActiveWorkbook

$template = $wb.Worksheets.Item("WFA (6)")
$playerSummary = $wb.Worksheets.Item("Player Summary")

# Copy WFA (6) so the new sheet lands immediately before Player Summary,
# bringing along its cell values/styles and conditional formatting. The
# new copy becomes the active sheet, so grab it from there rather than
# guessing Excel's default "WFA (6) (2)" name.
$template.Copy($playerSummary, $null)
$newSheet = $wb.ActiveSheet

# Rename it to follow the existing "WFA (n)" naming convention.
$newSheet.Name = "WFA (7)"

# Sheet.Copy() does not bring along the source ListObject/table, so
# recreate it on the new sheet over the same A3:I4 range, following the
# same incremental table-naming pattern used by the other WFA sheets
# (Table1, Table14, Table145, Table1456, Table14567, Table145678, ...).
$lo = $newSheet.ListObjects.Add(1, $newSheet.Range("A3:I4"), $null, 1)
$lo.Name = "Table1456789"
